$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh rolls the price history forward by one reporting
# period: a brand-new date/origin pair is inserted at the top of the
# "Acelga" data block (pushing the existing rows down by two, which also
# naturally reproduces the duplicated trailing pair seen at the bottom
# of the refreshed range), carrying the same Primera/Segunda stats that
# used to belong to the oldest pair.

$ws.Rows("129:130").Insert()

# New row 129 (Primera)
$ws.Range("A129").Value = 11
$ws.Range("B129").Value = "Vega Monumental Concepción"
$ws.Range("C129").Value = "Bíobío"
$ws.Range("D129").Value = 44582
$ws.Range("E129").Value = 8
$ws.Range("F129").Value = 100112009
$ws.Range("G129").Value = "Acelga"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 200
$ws.Range("K129").Value = 600
$ws.Range("L129").Value = 700
$ws.Range("M129").Value = 650
$ws.Range("N129").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O129").Value = "Región Metropolitana"
$ws.Range("P129").Value = 650
$ws.Range("Q129").Value = 1
$ws.Range("R129").Value = "Hortaliza"

# New row 130 (Segunda)
$ws.Range("A130").Value = 11
$ws.Range("B130").Value = "Vega Monumental Concepción"
$ws.Range("C130").Value = "Bíobío"
$ws.Range("D130").Value = 44582
$ws.Range("E130").Value = 8
$ws.Range("F130").Value = 100112009
$ws.Range("G130").Value = "Acelga"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Segunda"
$ws.Range("J130").Value = 100
$ws.Range("K130").Value = 500
$ws.Range("L130").Value = 500
$ws.Range("M130").Value = 500
$ws.Range("N130").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O130").Value = "Región Metropolitana"
$ws.Range("P130").Value = 500
$ws.Range("Q130").Value = 1
$ws.Range("R130").Value = "Hortaliza"
